$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51
# D-column values must remain plain text (e.g. "1.000", "29.210.15") so we
# force a text number format before assigning, then restore the default style.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.210.15'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.13%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.866.35'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.82%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7111'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.65%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.52'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.22%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07673'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.51%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '24.71'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.34%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08366'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.05%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.871.87'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.21%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.226'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.17%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7114'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.57%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.31'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.00%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.210.56'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.10%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.949'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.05%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '243.46'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.40%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007822'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.92%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.113.64'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.82%  '

$ws.Range('E21').Value = '  -2.05%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9992'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.13%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.858'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.49%  '

$ws.Range('E24').Value = '  -0.07%  '

$ws.Range('E25').Value = '  -1.37%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '163.26'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.02%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.955'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.26%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.49'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.71%  '

$ws.Range('E29').Value = '  -0.18%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.310'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.47%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.404'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.15%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.249'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.29%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05156'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.68%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7983'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +9.67%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.913'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.52%  '

$ws.Range('E36').Value = '  -2.87%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.684'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.06%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01854'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.80%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.712'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.14%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.159.21'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.10%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.307'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.57%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8953'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.62%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '73.16'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.08%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9997'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.13%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '103.08'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.80%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.010.69'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.84%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5184'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.90%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.780'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.97%  '

$ws.Range('E49').Value = '  -0.55%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.342'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.09%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4295'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.88%  '
